# "load data on web" / "load data from web in market_data lab"
# Adds a second worksheet (Sheet2) holding a "days / digoxin / change_digoxin"
# table, makes it the active sheet, and moves Sheet1's selection off J8.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Sheet1: selection moves from J8 to D2 (tabSelected will move to the new
# sheet once it is activated below).
[void]$ws1.Range("D2").Select()

# Insert the new sheet right after Sheet1.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)

# Header row.
$ws2.Range("A1").Value = "days"
$ws2.Range("B1").Value = "digoxin"
$ws2.Range("C1").Value = "change_digoxin"

# Data rows (day, digoxin remaining fraction, change vs. previous day).
# (Plain decimal literals -- the engine's own float formatting reproduces
# the scientific-notation forms used in the canonical XML on save.)
$data = @(
    @(0, 0.5,   -0.155),
    @(1, 0.345, -0.107),
    @(2, 0.238, -0.074),
    @(3, 0.164, -0.051),
    @(4, 0.113, -0.035),
    @(5, 0.078, -0.024),
    @(6, 0.054, -0.017),
    @(7, 0.037, -0.011),
    @(8, 0.026, $null)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 1).Value = $data[$i][0]
    $ws2.Cells.Item($row, 2).Value = $data[$i][1]
    if ($null -ne $data[$i][2]) {
        $ws2.Cells.Item($row, 3).Value = $data[$i][2]
    }
}

# New sheet ends up selected/active, matching tabSelected + activeTab.
[void]$ws2.Range("A1:C10").Select()
$ws2.Activate()
